$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Semester Information")

# Delete the extra semester row (row 3 - "Test 2")
$ws.Rows.Item(3).Delete()

# Update A2 to the placeholder semester name
$ws.Range("A2").Value = "This is a semester name that should get replaced"

# Clear old date-formatting on B2:C2 before writing new text
$ws.Range("B2:C2").ClearFormats()

# Force B2/C2 to be stored as literal text (not auto-converted to date serials)
$ws.Range("B2").Formula = "=""1/1/1111"""
$ws.Range("C2").Formula = "=""9/9/9999"""
$ws.Range("B2:C2").Copy()
$ws.Range("B2:C2").PasteSpecial(-4163)
